$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new value would otherwise be
# auto-parsed as a number by Excel, so they stay text (source data is all text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated coin data (values, swaps, percentages).
$ws.Range("D2").Value = "41.665.75"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "2.169.46"
$ws.Range("E3").Value = "  -3.16%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "238.76"
$ws.Range("E5").Value = "  -2.05%  "
$ws.Range("E6").Value = "  -3.26%  "
$ws.Range("D7").Value = "72.26"
$ws.Range("E7").Value = "  -3.06%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.579"
$ws.Range("E9").Value = "  -5.24%  "
$ws.Range("D10").Value = "40.08"
$ws.Range("E10").Value = "  -6.63%  "
$ws.Range("E11").Value = "  -5.80%  "
$ws.Range("D12").Value = "54.40"
$ws.Range("E12").Value = "  -4.80%  "
$ws.Range("E13").Value = "  -3.42%  "
$ws.Range("D14").Value = "6.73"
$ws.Range("E14").Value = "  -4.30%  "
$ws.Range("D15").Value = "2.494.80"
$ws.Range("E15").Value = "  -3.12%  "
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("D17").Value = "2.161.19"
$ws.Range("E17").Value = "  -2.96%  "
$ws.Range("D18").Value = "0.782"
$ws.Range("E18").Value = "  -7.29%  "
$ws.Range("D19").Value = "41.496.58"
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("D20").Value = "0.0000103"
$ws.Range("E20").Value = "  -2.85%  "
$ws.Range("D21").Value = "70.11"
$ws.Range("E21").Value = "  -4.28%  "
$ws.Range("E22").Value = "  -7.56%  "
$ws.Range("D23").Value = "9.82"
$ws.Range("E23").Value = "  -13.57%  "
$ws.Range("D24").Value = "226.98"
$ws.Range("E24").Value = "  -2.11%  "
$ws.Range("D25").Value = "2.00"
$ws.Range("E25").Value = "  -5.07%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -6.41%  "
$ws.Range("D28").Value = "3.27"
$ws.Range("E28").Value = "  -9.77%  "
$ws.Range("E29").Value = "  -4.06%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.16"
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "170.56"
$ws.Range("E31").Value = "  +1.99%  "
$ws.Range("E32").Value = "  -4.14%  "
$ws.Range("D33").Value = "33.44"
$ws.Range("E33").Value = "  +8.89%  "
$ws.Range("D34").Value = "0.0774"
$ws.Range("E34").Value = "  -4.05%  "
$ws.Range("D35").Value = "5.22"
$ws.Range("E35").Value = "  -8.76%  "
$ws.Range("E36").Value = "  -3.61%  "
$ws.Range("D37").Value = "4.33"
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("E38").Value = "  -4.06%  "
$ws.Range("D39").Value = "0.0309"
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "12.11"
$ws.Range("E40").Value = "  -10.51%  "
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").Value = "2.10"
$ws.Range("E41").Value = "  -2.42%  "
$ws.Range("E42").Value = "  -6.12%  "
$ws.Range("D43").Value = "59.06"
$ws.Range("E43").Value = "  -9.19%  "
$ws.Range("E44").Value = "  -3.76%  "
$ws.Range("E45").Value = "  -5.76%  "
$ws.Range("E46").Value = "  -4.36%  "
$ws.Range("D47").Value = "97.59"
$ws.Range("E47").Value = "  -7.36%  "
$ws.Range("E48").Value = "  -4.31%  "
$ws.Range("D49").Value = "1.12"
$ws.Range("E49").Value = "  -4.93%  "
$ws.Range("D50").Value = "2.18"
$ws.Range("E50").Value = "  -8.01%  "
$ws.Range("E51").Value = "  -2.41%  "
